$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Timer Resolution" (B6) changes from 1E-3 seconds to 1E-5 seconds
$ws.Range("B6").Value = 0.00001

# Recalculate so the dependent formula in B8 (=B6/B3) updates its cached value
$excel.Calculate()

# Move the active selection from B8 to B7
$ws.Range("B7").Select()
